$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: Delete column BK entirely.
# Old columns BL:BQ (header text in row 1, values in rows 2-85) shift left into BK:BP,
# matching the relabeling/reordering seen in the diff.
$ws.Range("BK:BK").Delete()

# Step 2: Write headers for the now-empty BQ1 (re-added "night avg temp" label)
# and the two brand new columns BR1, BS1.
$ws.Range("BQ1").Value = "야간평균온도"
$ws.Range("BR1").Value = "일몰일출적합증산(HD)누적시간"
$ws.Range("BS1").Value = "주야간온도차이"

# Copy the header formatting (bold, centered, thin border) from an existing header
# cell onto the newly written header cells so they match the rest of row 1.
$ws.Range("BP1").Copy()
$ws.Range("BQ1:BS1").PasteSpecial(-4122)

# Step 3: Write the recalculated data values for BQ, BR, BS across data rows 2-84.
$ws.Range("BQ2").Value = 17.51379362670719
$ws.Range("BR2").Value = 7.916666666666667
$ws.Range("BS2").Value = 4.218920841922877
$ws.Range("BQ3").Value = 14.17154779969653
$ws.Range("BR3").Value = 0
$ws.Range("BS3").Value = 3.264529123380502
$ws.Range("BQ4").Value = 14.19098784194526
$ws.Range("BR4").Value = 0
$ws.Range("BS4").Value = 5.18722644376915
$ws.Range("BQ5").Value = 15.23522003034902
$ws.Range("BR5").Value = 3.25
$ws.Range("BS5").Value = 3.526446636317759
$ws.Range("BQ6").Value = 14.43637329286799
$ws.Range("BR6").Value = 3.616666666666667
$ws.Range("BS6").Value = 3.694442769308282
$ws.Range("BQ7").Value = 14.48665130568356
$ws.Range("BR7").Value = 0
$ws.Range("BS7").Value = 4.840493380616202
$ws.Range("BQ8").Value = 14.42584218512891
$ws.Range("BR8").Value = 2.466666666666667
$ws.Range("BS8").Value = 4.01543822460231
$ws.Range("BQ9").Value = 14.77216923076928
$ws.Range("BR9").Value = 0
$ws.Range("BS9").Value = 4.218694171292684
$ws.Range("BQ10").Value = 14.85311248073963
$ws.Range("BR10").Value = 0.5333333333333333
$ws.Range("BS10").Value = 4.415892687219129
$ws.Range("BQ11").Value = 16.28509345794391
$ws.Range("BR11").Value = 2.616666666666667
$ws.Range("BS11").Value = 4.84788831288952
$ws.Range("BQ12").Value = 15.45644376899704
$ws.Range("BR12").Value = 5.6
$ws.Range("BS12").Value = 6.896422638266367
$ws.Range("BQ13").Value = 16.20408194233684
$ws.Range("BR13").Value = 0
$ws.Range("BS13").Value = 1.611610365355631
$ws.Range("BQ14").Value = 17.13142640364192
$ws.Range("BR14").Value = 0
$ws.Range("BS14").Value = 1.634586416871006
$ws.Range("BQ15").Value = 15.88248861911995
$ws.Range("BR15").Value = 7.966666666666667
$ws.Range("BS15").Value = 4.824143903287393
$ws.Range("BQ16").Value = 15.83864946889233
$ws.Range("BR16").Value = 0.8166666666666667
$ws.Range("BS16").Value = 0.623693680915661
$ws.Range("BQ17").Value = 14.41911987860393
$ws.Range("BR17").Value = 4.016666666666667
$ws.Range("BS17").Value = 2.230483194379445
$ws.Range("BQ18").Value = 13.96917933130703
$ws.Range("BR18").Value = 10.96666666666667
$ws.Range("BS18").Value = 3.171090245073035
$ws.Range("BQ19").Value = 12.86922492401221
$ws.Range("BR19").Value = 5.85
$ws.Range("BS19").Value = 1.023082768295435
$ws.Range("BQ20").Value = 13.8677996965099
$ws.Range("BR20").Value = 10.98333333333333
$ws.Range("BS20").Value = 3.802481993631005
$ws.Range("BQ21").Value = 13.84667173252282
$ws.Range("BR21").Value = 10.96666666666667
$ws.Range("BS21").Value = 2.343674865680061
$ws.Range("BQ22").Value = 15.45205167173263
$ws.Range("BR22").Value = 10.96666666666667
$ws.Range("BS22").Value = 1.970227457588841
$ws.Range("BQ23").Value = 16.25084977238242
$ws.Range("BR23").Value = 10.38333333333333
$ws.Range("BS23").Value = 2.694361495223351
$ws.Range("BQ24").Value = 15.28931714719273
$ws.Range("BR24").Value = 10.98333333333333
$ws.Range("BS24").Value = 3.966867231808711
$ws.Range("BQ25").Value = 15.27095599393018
$ws.Range("BR25").Value = 10.98333333333333
$ws.Range("BS25").Value = 4.654101624507863
$ws.Range("BQ26").Value = 15.41637329286812
$ws.Range("BR26").Value = 10.98333333333333
$ws.Range("BS26").Value = 5.130969454243681
$ws.Range("BQ27").Value = 14.65984825493179
$ws.Range("BR27").Value = 0.3833333333333334
$ws.Range("BS27").Value = 5.593728668145234
$ws.Range("BQ28").Value = 13.64753424657538
$ws.Range("BR28").Value = 5.333333333333333
$ws.Range("BS28").Value = 8.011134127304393
$ws.Range("BQ29").Value = 13.7506775700935
$ws.Range("BR29").Value = 5.916666666666667
$ws.Range("BS29").Value = 7.647078840163042
$ws.Range("BQ30").Value = 14.09015220700154
$ws.Range("BR30").Value = 0.2666666666666667
$ws.Range("BS30").Value = 5.669988818639567
$ws.Range("BQ31").Value = 14.20526555386952
$ws.Range("BR31").Value = 1.816666666666667
$ws.Range("BS31").Value = 8.929016497412649
$ws.Range("BQ32").Value = 16.98000000000011
$ws.Range("BR32").Value = 0.9
$ws.Range("BS32").Value = 6.588202824133536
$ws.Range("BQ33").Value = 15.46848024316111
$ws.Range("BR33").Value = 0
$ws.Range("BS33").Value = 1.632390435456157
$ws.Range("BQ34").Value = 13.95250379362678
$ws.Range("BR34").Value = 0
$ws.Range("BS34").Value = 0.7721953100864347
$ws.Range("BQ35").Value = 14.32415781487108
$ws.Range("BR35").Value = 5.25
$ws.Range("BS35").Value = 3.075739752350504
$ws.Range("BQ36").Value = 13.87804281345568
$ws.Range("BR36").Value = 10.9
$ws.Range("BS36").Value = 3.606231897712593
$ws.Range("BQ37").Value = 15.97301972685888
$ws.Range("BR37").Value = 10.98333333333333
$ws.Range("BS37").Value = 3.908252766714522
$ws.Range("BQ38").Value = 14.32197268588775
$ws.Range("BR38").Value = 10.98333333333333
$ws.Range("BS38").Value = 5.376849337159745
$ws.Range("BQ39").Value = 18.1224734446131
$ws.Range("BR39").Value = 3.283333333333333
$ws.Range("BS39").Value = 5.078309610585958
$ws.Range("BQ40").Value = 15.24474962063734
$ws.Range("BR40").Value = 8.733333333333333
$ws.Range("BS40").Value = 6.463419393447303
$ws.Range("BQ41").Value = 13.76599383667186
$ws.Range("BR41").Value = 9.483333333333333
$ws.Range("BS41").Value = 7.383314742073475
$ws.Range("BQ42").Value = 14.56762658227856
$ws.Range("BR42").Value = 2.45
$ws.Range("BS42").Value = 4.488517376590446
$ws.Range("BQ43").Value = 14.14383915022765
$ws.Range("BR43").Value = 3.15
$ws.Range("BS43").Value = 0.8030711061826956
$ws.Range("BQ44").Value = 14.39303951367785
$ws.Range("BR44").Value = 10.96666666666667
$ws.Range("BS44").Value = 3.915484234717599
$ws.Range("BQ45").Value = 13.40549317147196
$ws.Range("BR45").Value = 9.016666666666667
$ws.Range("BS45").Value = 5.802122213143601
$ws.Range("BQ46").Value = 14.09295558958655
$ws.Range("BR46").Value = 9.35
$ws.Range("BS46").Value = 6.870977572367302
$ws.Range("BQ47").Value = 14.48669195751141
$ws.Range("BR47").Value = 6.25
$ws.Range("BS47").Value = 6.150509582925208
$ws.Range("BQ48").Value = 13.98106221547805
$ws.Range("BR48").Value = 10.01666666666667
$ws.Range("BS48").Value = 7.687848040932357
$ws.Range("BQ49").Value = 15.10103343465049
$ws.Range("BR49").Value = 8.4
$ws.Range("BS49").Value = 10.08355043218702
$ws.Range("BQ50").Value = 15.05001517450681
$ws.Range("BR50").Value = 1.55
$ws.Range("BS50").Value = 9.422869440877912
$ws.Range("BQ51").Value = 17.25875379939212
$ws.Range("BR51").Value = 1.883333333333333
$ws.Range("BS51").Value = 3.139536689039868
$ws.Range("BQ52").Value = 14.34077389984826
$ws.Range("BR52").Value = 0
$ws.Range("BS52").Value = 1.756870146246541
$ws.Range("BQ53").Value = 14.26418816388472
$ws.Range("BR53").Value = 0
$ws.Range("BS53").Value = 4.028337477141015
$ws.Range("BQ54").Value = 14.55830534351151
$ws.Range("BR54").Value = 4.466666666666667
$ws.Range("BS54").Value = 7.478310041104038
$ws.Range("BQ55").Value = 14.79068285280732
$ws.Range("BR55").Value = 8.833333333333334
$ws.Range("BS55").Value = 6.784082772192807
$ws.Range("BQ56").Value = 14.26689969604866
$ws.Range("BR56").Value = 1.783333333333333
$ws.Range("BS56").Value = 7.127095182312578
$ws.Range("BQ57").Value = 15.47995447647953
$ws.Range("BR57").Value = 10.53333333333333
$ws.Range("BS57").Value = 6.451121067694753
$ws.Range("BQ58").Value = 15.0179969650986
$ws.Range("BR58").Value = 1.05
$ws.Range("BS58").Value = 7.75400047408205
$ws.Range("BQ59").Value = 18.83951367781164
$ws.Range("BR59").Value = 1
$ws.Range("BS59").Value = 5.945825630767199
$ws.Range("BQ60").Value = 20.27308980213105
$ws.Range("BR60").Value = 0.5333333333333333
$ws.Range("BS60").Value = 4.70069710415919
$ws.Range("BQ61").Value = 14.77113808801217
$ws.Range("BR61").Value = 0
$ws.Range("BS61").Value = 7.106211463844522
$ws.Range("BQ62").Value = 13.96793626707135
$ws.Range("BR62").Value = 3.116666666666667
$ws.Range("BS62").Value = 5.349528521661169
$ws.Range("BQ63").Value = 14.50423368740519
$ws.Range("BR63").Value = 5.866666666666666
$ws.Range("BS63").Value = 7.186124827319686
$ws.Range("BQ64").Value = 18.61624620060811
$ws.Range("BR64").Value = 4.4
$ws.Range("BS64").Value = 5.070463146383041
$ws.Range("BQ65").Value = 17.71966616084976
$ws.Range("BR65").Value = 0
$ws.Range("BS65").Value = 1.357047574708748
$ws.Range("BQ66").Value = 15.64107902735571
$ws.Range("BR66").Value = 0
$ws.Range("BS66").Value = 5.715407459130914
$ws.Range("BQ67").Value = 14.04132218844986
$ws.Range("BR67").Value = 3.316666666666667
$ws.Range("BS67").Value = 8.482135876066389
$ws.Range("BQ68").Value = 13.90861911987864
$ws.Range("BR68").Value = 0.8833333333333333
$ws.Range("BS68").Value = 8.093314298815512
$ws.Range("BQ69").Value = 14.93772382397565
$ws.Range("BR69").Value = 0.8166666666666667
$ws.Range("BS69").Value = 6.569817789340782
$ws.Range("BQ70").Value = 16.73487101669205
$ws.Range("BR70").Value = 3.683333333333333
$ws.Range("BS70").Value = 5.194628341459548
$ws.Range("BQ71").Value = 17.30430091185421
$ws.Range("BR71").Value = 0
$ws.Range("BS71").Value = 1.471064004919235
$ws.Range("BQ72").Value = 14.75894817073171
$ws.Range("BR72").Value = 0
$ws.Range("BS72").Value = 6.814368090472023
$ws.Range("BQ73").Value = 17.3675379939209
$ws.Range("BR73").Value = 0.7666666666666667
$ws.Range("BS73").Value = 6.090618216066442
$ws.Range("BQ74").Value = 23.35389057750776
$ws.Range("BR74").Value = 0
$ws.Range("BS74").Value = 2.916506858389788
$ws.Range("BQ75").Value = 24.33557926829291
$ws.Range("BR75").Value = 0
$ws.Range("BS75").Value = 2.006699861028562
$ws.Range("BQ76").Value = 22.68773899848278
$ws.Range("BR76").Value = 0
$ws.Range("BS76").Value = 4.561889682695288
$ws.Range("BQ77").Value = 17.11983257229833
$ws.Range("BR77").Value = 0
$ws.Range("BS77").Value = 4.956842652348847
$ws.Range("BQ78").Value = 15.28905775075987
$ws.Range("BR78").Value = 9.733333333333333
$ws.Range("BS78").Value = 7.109700251801076
$ws.Range("BQ79").Value = 14.53301972685887
$ws.Range("BR79").Value = 1.166666666666667
$ws.Range("BS79").Value = 8.543895440236385
$ws.Range("BQ80").Value = 16.60655538694993
$ws.Range("BR80").Value = 1.3
$ws.Range("BS80").Value = 7.486184689874737
$ws.Range("BQ81").Value = 19.27007587253434
$ws.Range("BR81").Value = 0
$ws.Range("BS81").Value = 4.43946317996258
$ws.Range("BQ82").Value = 13.99364741641346
$ws.Range("BR82").Value = 1.166666666666667
$ws.Range("BS82").Value = 11.17860610471339
$ws.Range("BQ83").Value = 14.57723404255319
$ws.Range("BR83").Value = 1.2
$ws.Range("BS83").Value = 9.534612111293104
$ws.Range("BQ84").Value = 16.43981790591802
$ws.Range("BR84").Value = 1.3
$ws.Range("BS84").Value = 8.253041068441089

# Step 4: Remove the old row 85, which no longer exists in the updated dataset.
$ws.Range("85:85").Delete()
